$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they remain text (matching original inlineStr type)
$numericLikeCells = @("D5", "D6", "D8", "D10", "D12", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D51")
foreach ($cellRef in $numericLikeCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "60.824.68"
$ws.Range("E2").Value = "  +3.80%  "
$ws.Range("D3").Value = "2.715.19"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "528.07"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "146.90"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "2.739.76"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  +14.23%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "3.198.21"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "60.823.71"
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").Value = "2.796.56"
$ws.Range("E17").Value = "  +6.36%  "
$ws.Range("D18").Value = "0.0000138"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "345.67"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "4.52"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "10.59"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  +5.09%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "63.26"
$ws.Range("E24").Value = "  +2.69%  "
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  +5.19%  "
$ws.Range("D26").Value = "0.418"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "0.995"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "0.0₃0825"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  +4.94%  "
$ws.Range("D30").Value = "6.76"
$ws.Range("E30").Value = "  +9.14%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "19.10"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").Value = "150.25"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +7.69%  "
$ws.Range("D36").Value = "1.23"
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("D37").Value = "0.921"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").Value = "0.901"
$ws.Range("E38").Value = "  +8.21%  "
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  +8.31%  "
$ws.Range("D40").Value = "37.24"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "3.69"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "0.627"
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("D43").Value = "20.38"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "281.79"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "0.0989"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "2.117.04"
$ws.Range("E47").Value = "  +6.90%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "4.94"
$ws.Range("E48").Value = "  +6.88%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0543"
$ws.Range("E49").Value = "  +4.32%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "19.37"
$ws.Range("E51").Value = "  +5.72%  "
